# Results of LSTM with wind speed
# Update the "LSTM, past p 168h, ws 1h forecast" row (row 7) confidence
# interval values, and move the active cell selection to E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C7").Value = "±2.96"
$ws.Range("D7").Value = "±4.22"
$ws.Range("E7").Value = "±4.91"

$ws.Range("E8").Select()
